# Rotates the data of rows 3, 4 and 5 (A32615-2023 observation records):
#   new row 3 <- old row 5
#   new row 4 <- old row 3
#   new row 5 <- old row 4
# Only the columns that actually differ between the rotated rows are
# touched (A, B, D, E, F, G, H, I, Q, R, S, Z, AB) - everything else
# (C, T, U, V, W, Y, AA, AD, AE, AG, AT, AW, AX, AY, ...) is identical
# across the three rows already, so it is left alone.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- snapshot the "before" values we need, straight from the sheet -------
$row3 = @{
    A = $ws.Range("A3").Value2
    B = $ws.Range("B3").Value2
    D = $ws.Range("D3").Value2
    E = $ws.Range("E3").Value2
    F = $ws.Range("F3").Value2
    G = $ws.Range("G3").Value2
    H = $ws.Range("H3").Value2
    I = $ws.Range("I3").Value2
    Q = $ws.Range("Q3").Value2
    R = $ws.Range("R3").Value2
    S = $ws.Range("S3").Value2
    Z = $ws.Range("Z3").Value2
    AB = $ws.Range("AB3").Value2
}
$row4 = @{
    A = $ws.Range("A4").Value2
    B = $ws.Range("B4").Value2
    D = $ws.Range("D4").Value2
    E = $ws.Range("E4").Value2
    F = $ws.Range("F4").Value2
    G = $ws.Range("G4").Value2
    H = $ws.Range("H4").Value2
    I = $ws.Range("I4").Value2
    Q = $ws.Range("Q4").Value2
    R = $ws.Range("R4").Value2
    S = $ws.Range("S4").Value2
    Z = $ws.Range("Z4").Value2
    AB = $ws.Range("AB4").Value2
}
$row5 = @{
    A = $ws.Range("A5").Value2
    B = $ws.Range("B5").Value2
    D = $ws.Range("D5").Value2
    E = $ws.Range("E5").Value2
    F = $ws.Range("F5").Value2
    G = $ws.Range("G5").Value2
    H = $ws.Range("H5").Value2
    I = $ws.Range("I5").Value2
    Q = $ws.Range("Q5").Value2
    R = $ws.Range("R5").Value2
    S = $ws.Range("S5").Value2
    Z = $ws.Range("Z5").Value2
    AB = $ws.Range("AB5").Value2
}

function Set-ArtfyndRow($destRow, $src) {
    $ws.Range("A$destRow").Value = $src.A
    $ws.Range("B$destRow").Value = $src.B
    $ws.Range("D$destRow").Value = $src.D
    $ws.Range("E$destRow").Value = $src.E
    $ws.Range("F$destRow").Value = $src.F
    $ws.Range("G$destRow").Value = $src.G
    $ws.Range("H$destRow").Value = $src.H

    # "Antal" (I) is stored as text in the source data (e.g. "35"), even
    # though it looks numeric, so force text formatting before assigning
    # it - otherwise Excel would auto-detect it as a number.
    if ($null -eq $src.I -or $src.I -eq "") {
        $ws.Range("I$destRow").Value = ""
    } else {
        $ws.Range("I$destRow").NumberFormat = "@"
        $ws.Range("I$destRow").Value = [string]$src.I
    }

    $ws.Range("Q$destRow").Value = $src.Q
    $ws.Range("R$destRow").Value = $src.R
    $ws.Range("S$destRow").Value = $src.S
    $ws.Range("Z$destRow").Value = $src.Z
    $ws.Range("AB$destRow").Value = $src.AB
}

# --- write the rotated rows back ------------------------------------------
Set-ArtfyndRow 3 $row5
Set-ArtfyndRow 4 $row3
Set-ArtfyndRow 5 $row4

Write-Output "rotation applied"
